$d = $word.ActiveDocument

# --- 1) "o)" paragraph: color the existing text (o, ) , Haz que..., Seguro que
#        quieres salir?) in red (FF0000). The surrounding <w:br/> runs and the
#        proofErr markers stay untouched. ---
$rngO = $d.Content
$targetO = "o) Haz que cuando se pulse el botón atrás físico del teléfono que se muestre un cuadro de mensaje que diga. Seguro que quieres salir?"
$rngO.Find.Execute($targetO, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rngO.Font.Color = 255

# --- 2) "t)" paragraph: color the whole exercise text red. ---
$rngT = $d.Content
$targetT = "t) Crea una clase deportista  con los valores de: puntuación, nombre, deporte y edad. Estos deben mostrarse en una grid. "
$rngT.Find.Execute($targetT, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rngT.Font.Color = 255

# --- 3) second "u)" paragraph (mapa exercise): color the whole exercise text
#        red, then move the _GoBack bookmark to right after it (adding a
#        bookmark with an existing name relocates it, removing the old one
#        that used to sit after the "picker" sentence). ---
$rngU = $d.Content
$targetU = "u) Crea un mapa y muestra la ubicación del estadio del Betis o el Sevilla y la ubicación  del Fesac. En el mismo."
$rngU.Find.Execute($targetU, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rngU.Font.Color = 255
$rngU.Collapse(0)
$d.Bookmarks.Add("_GoBack", $rngU)
